$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns (AD1:AF1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (e.g. from AC1) to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

for ($row = 2; $row -le 56; $row++) {
    $ws.Cells.Item($row, 30).Value = 68   # AD
    $ws.Cells.Item($row, 31).Value = 94   # AE
    $ws.Cells.Item($row, 32).Value = 0    # AF
}
